$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.128.53"
$ws.Range("E2").Value = "  -0.27%  "
$ws.Range("D3").Value = "2.500.83"
$ws.Range("E3").Value = "  -0.38%  "
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "'537.29"
$ws.Range("E5").Value = "  -0.41%  "
$ws.Range("D6").Value = "'137.07"
$ws.Range("E6").Value = "  -1.41%  "
$ws.Range("D7").Value = "'0.997"
$ws.Range("E7").Value = "  -0.23%  "
$ws.Range("E8").Value = "  +0.79%  "
$ws.Range("D9").Value = "2.525.29"
$ws.Range("E9").Value = "  +0.55%  "
$ws.Range("E10").Value = "  -0.09%  "
$ws.Range("E11").Value = "  -2.10%  "
$ws.Range("D12").Value = "'5.33"
$ws.Range("E12").Value = "  -1.01%  "
$ws.Range("E13").Value = "  -1.28%  "
$ws.Range("D14").Value = "2.948.55"
$ws.Range("E14").Value = "  -0.50%  "
$ws.Range("E15").Value = "  -1.07%  "
$ws.Range("D16").Value = "58.934.39"
$ws.Range("E16").Value = "  -0.46%  "
$ws.Range("E17").Value = "  -0.96%  "
$ws.Range("D18").Value = "2.520.78"
$ws.Range("E18").Value = "  +0.42%  "
$ws.Range("D19").Value = "'11.15"
$ws.Range("E19").Value = "  +0.79%  "
$ws.Range("E20").Value = "  +0.05%  "
$ws.Range("D21").Value = "'324.24"
$ws.Range("E21").Value = "  -0.17%  "
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("E23").Value = "  +1.50%  "
$ws.Range("D24").Value = "'65.80"
$ws.Range("E24").Value = "  +3.77%  "
$ws.Range("E25").Value = "  +0.19%  "
$ws.Range("E26").Value = "  -1.81%  "
$ws.Range("D27").Value = "'0.997"
$ws.Range("E27").Value = "  -0.19%  "
$ws.Range("E28").Value = "  -3.29%  "
$ws.Range("D29").Value = "'6.72"
$ws.Range("E29").Value = "  -1.32%  "
$ws.Range("D30").Value = "0.0₃0775"
$ws.Range("E30").Value = "  -0.15%  "
$ws.Range("E31").Value = "  -1.20%  "
$ws.Range("E32").Value = "  +2.61%  "
$ws.Range("E33").Value = "  +6.20%  "
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("E35").Value = "  +1.84%  "
$ws.Range("D37").Value = "'4.12"
$ws.Range("E37").Value = "  -2.30%  "
$ws.Range("E38").Value = "  -2.95%  "
$ws.Range("D39").Value = "'36.71"
$ws.Range("E39").Value = "  -0.58%  "
$ws.Range("D40").Value = "'0.816"
$ws.Range("E40").Value = "  +0.73%  "
$ws.Range("E41").Value = "  -1.29%  "
$ws.Range("D42").Value = "'284.82"
$ws.Range("E42").Value = "  +1.74%  "
$ws.Range("D43").Value = "'5.23"
$ws.Range("E43").Value = "  +0.34%  "
$ws.Range("D44").Value = "'131.91"
$ws.Range("E44").Value = "  +5.30%  "
$ws.Range("D45").Value = "'0.996"
$ws.Range("E45").Value = "  -0.23%  "
$ws.Range("E46").Value = "  +1.67%  "
$ws.Range("D47").Value = "'10.89"
$ws.Range("E47").Value = "  +0.05%  "
$ws.Range("E48").Value = "  -1.11%  "
$ws.Range("E49").Value = "  -0.79%  "
$ws.Range("E50").Value = "  -1.64%  "
$ws.Range("D51").Value = "'17.37"
$ws.Range("E51").Value = "  -2.36%  "
